$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 239; existing rows 239:252 shift down to 240:253.
$ws.Rows(239).Insert()

$ws.Cells.Item(239, 1).Value = 3
$ws.Cells.Item(239, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(239, 3).Value = "Coquimbo"
$ws.Cells.Item(239, 4).Value = 44516
$ws.Cells.Item(239, 5).Value = 5
$ws.Cells.Item(239, 6).Value = 100112003
$ws.Cells.Item(239, 7).Value = "Ajo"
$ws.Cells.Item(239, 8).Value = "Chino"
$ws.Cells.Item(239, 9).Value = "Primera"
$ws.Cells.Item(239, 10).Value = 75
$ws.Cells.Item(239, 11).Value = 16000
$ws.Cells.Item(239, 12).Value = 16500
$ws.Cells.Item(239, 13).Value = 16267
$ws.Cells.Item(239, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(239, 15).Value = "China"
$ws.Cells.Item(239, 16).Value = 1627
$ws.Cells.Item(239, 17).Value = 10
$ws.Cells.Item(239, 18).Value = "Hortaliza"
